$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = 13.191
$ws.Range("B12").Value = 6.112
$ws.Range("E12").Value = 13.143
$ws.Range("E14").Value = 13.072
$ws.Range("E22").Value = 13.126
